$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33
$ws.Cells.Item(33, 1).Value = "Anil Prasad"
$ws.Cells.Item(33, 2).Value = "kumarr47872@gmail.com"
$ws.Cells.Item(33, 3).Value = "Paris"
$ws.Cells.Item(33, 4).Value = "'2025-02-18"
$ws.Cells.Item(33, 5).Value = "'4"
$ws.Cells.Item(33, 6).Value = "Nothing"
$ws.Cells.Item(33, 7).Value = "'2025-02-17 23:05:28"

# Row 34
$ws.Cells.Item(34, 1).Value = "RAHUL KUMAR"
$ws.Cells.Item(34, 2).Value = "231fa04862@gmail.com"
$ws.Cells.Item(34, 3).Value = "Paris"
$ws.Cells.Item(34, 4).Value = "'2025-02-19"
$ws.Cells.Item(34, 5).Value = "'4"
$ws.Cells.Item(34, 6).Value = "Nothing "
$ws.Cells.Item(34, 7).Value = "'2025-02-17 23:12:05"
